$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap header labels: A1 <-> B1 (Tank_ID / time)
$ws.Range("A1").Value = "Tank_ID"
$ws.Range("B1").Value = "time"

# Swap data values: A2 <-> B2 (1 / 7200)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 7200

# Update the active selection to C14
$ws.Range("C14").Select()
